$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header in B1: "Liked" -> "judgement"
$ws.Range("B1").Value = "judgement"

# Replace numeric 1/0 values in B2:B36 with text labels "good"/"bad",
# and right-align those cells.
$labels = @("good","bad","bad","good","bad","good","good","bad","bad","good","bad",
            "good","bad","bad","good","good","good","bad","bad","bad","good","bad",
            "bad","good","good","good","good","good","bad","bad","bad","good","good",
            "bad","good")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $labels[$i]
    $cell.HorizontalAlignment = -4152
}

# Update the sheet view: scroll so row 1 is visible again and move the
# active selection to F7 (matches the saved view state in the workbook).
$ws.Application.GoTo($ws.Range("F7"))
